# Commit: "delete teaching materials and change the corresponding description"
#
# The sheet has a list of repeated "section" blocks (title row + header row +
# a handful of blank data rows). The old layout had two back-to-back
# sections:
#   row 30      : "课件"   (section title)
#   row 31      : 课程名 / 文件夹 / 授课老师 / 文件个数 / 最近更新时间  (headers)
#   rows 32-36  : blank data rows
#   row 37      : "教材"   (section title)
#   row 38      : 课程名 / 教材名 / 作者 / 出版时间 / 出版社 / ISBN   (headers)
#   rows 39-43  : blank data rows
#
# The edit removes the old "课件" section entirely (title + header + its 5
# blank rows) so the "教材" section shifts up into its place, becoming the
# new rows 30-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "课件" section (rows 30-36): title row, header row, and its
# 5 blank data rows. Everything below (the "教材" section) shifts up to
# fill the gap.
$ws.Range("A30:A36").EntireRow.Delete()

# Match the author's saved cursor position.
$ws.Range("J32").Select()
